# Update the "Gaz" worksheet with the new EPEX Spot price columns
# (Last Price, Last Volume, End of Day Index) as produced by the
# automated data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# --- Header row (row 1): new column titles ---
$ws.Range("E1").Value = "Last Price"
$ws.Range("F1").Value = "Last Volume"
$ws.Range("G1").Value = "End of Day Index"

# Match the formatting already used by the existing header cells
# (bold font, thin border, centered alignment) by copying D1's style.
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 2: no data available yet for the new columns ---
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# --- Row 3: latest EPEX spot values ---
$ws.Range("E3").Value = 38.95
$ws.Range("F3").Value = 24000
$ws.Range("G3").Value = 38.201
